$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Annotation matrix values for columns E:J, rows 2-17
$data = @(
    @(2,2,1,2,1,2),
    @(2,2,1,2,1,2),
    @(2,1,1,1,1,2),
    @(2,0,1,1,2,2),
    @(2,2,2,2,1,2),
    @(2,1,1,1,1,2),
    @(2,0,1,2,1,2),
    @(2,1,2,2,1,2),
    @(2,2,1,2,2,2),
    @(2,2,2,1,1,2),
    @(2,1,1,0,1,2),
    @(2,2,2,2,1,2),
    @(2,2,1,1,1,2),
    @(2,2,1,1,1,2),
    @(2,2,1,1,1,2),
    @(2,1,1,1,2,2)
)

$startRow = 2
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $vals = $data[$i]
    for ($j = 0; $j -lt $vals.Length; $j++) {
        $col = 5 + $j  # E=5
        $ws.Cells.Item($row, $col).Value = $vals[$j]
    }
}

# Set zoom and freeze panes on the sheet view
$ws.Activate()
$excel.ActiveWindow.Zoom = 85

$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

$ws.Range("A8").Select()
$ws.Range("E18").Select()
